# Automatic update of files.
# Swap the content of row 3 and row 4 (the two observation records got
# reordered / renumbered), moving each record's data to the other row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: Id ---
$ws.Range("A3").Value = 111697636
$ws.Range("A4").Value = 111697236

# --- Column B: Taxonsorteringsordning ---
$ws.Range("B3").Value = 88489
$ws.Range("B4").Value = 8377

# --- Column D: Rödlistade ---
$ws.Range("D3").Value = "NT"
$ws.Range("D4").Value = "LC"

# --- Column E: TaxonId ---
$ws.Range("E3").Value = 1962
$ws.Range("E4").Value = 106545

# --- Column F: Artnamn ---
$ws.Range("F3").Value = "Vaddporing"
$ws.Range("F4").Value = "Mindre märgborre"

# --- Column G: Vetenskapligt namn ---
$ws.Range("G3").Value = "Anomoporia kamtschatica"
$ws.Range("G4").Value = "Tomicus minor"

# --- Column H: Auktor ---
$ws.Range("H3").Value = "(Parmasto) Bondartseva"
$ws.Range("H4").Value = "(Hartig, 1834)"

# --- Column J: Enhet (row3 gains "fruktkroppar", row4 loses it) ---
$ws.Range("J3").Value = "fruktkroppar"
$ws.Range("J4").Value = ""

# --- Columns L/M: Kön / Aktivitet - move from row3 to row4 ---
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = "färska gnagspår"

# --- Column Q: Ost ---
$ws.Range("Q3").Value = 373112.5181173298
$ws.Range("Q4").Value = 373121.3523494597

# --- Column R: Nord ---
$ws.Range("R3").Value = 6865358.590016441
$ws.Range("R4").Value = 6865443.651501717

# --- Column Z: Starttid ---
$ws.Range("Z3").Value = "19:00"
$ws.Range("Z4").Value = "00:00"

# --- Column AB: Sluttid ---
$ws.Range("AB3").Value = "19:00"
$ws.Range("AB4").Value = "00:00"

# --- Column AC: Publik kommentar - move from row4 to row3 ---
$ws.Range("AC3").Value = "Växer under rötad gammal silverved"
$ws.Range("AC4").ClearContents()

# --- Column AI: Biotop-beskrivning ---
$ws.Range("AI3").Value = "Kontinuitetsskog. Tallskog"
$ws.Range("AI4").Value = "Tallskog. Kontinuitetsskog"
